$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly by writing values; write rows 2..31 for columns A-H
$data = @(
    @(0, "struggle", 0.037449836730957, 0.3890565633773803, -0.2103400230407714, 0.00580321977447178, 0.07317293116024544, 0.01743147575429496),
    @(100, "struggle", 0.0105371475219726, 0.3668201565742492, -0.3450851440429687, -0.01212380492907687, 0.1479228914392232, 0.03119464347861243),
    @(200, "struggle", 0.1156568527221679, 0.2621434330940246, 0.2320724725723266, -0.05563860289676467, 0.2476778115545005, 0.03135671048444155),
    @(300, "struggle", -0.012700080871582, 0.3827533721923828, 0.0918664336204528, 0.000710598117082524, -0.0459115079471038, -0.03720979673825954),
    @(400, "struggle", -0.1776895523071289, 1.103561997413635, -1.161585211753845, 0.03395912094915998, -0.002209710983597735, -0.02471199281969851),
    @(500, "struggle", -0.1967945098876953, -0.0881298780441284, -0.4354097247123718, 0.0100699262883591, -0.003861541194575109, -0.02743906991518268),
    @(600, "struggle", 0.3726930618286133, 0.3928739428520202, -0.1955753564834594, -0.04216528505238931, -0.05587235412427344, -0.005946585338334982),
    @(700, "struggle", 0.1368236541748047, 0.4095092415809631, -0.2153286337852478, -0.2539235970803638, -0.6386596262454889, 0.1363320350646959),
    @(800, "struggle", -0.0151729583740234, 0.2465144991874694, 0.0501556992530822, -1.24240333419673, -3.812540969678324, 0.5514710630689335),
    @(900, "struggle", 2.021584510803223, -1.973333358764648, -6.84404182434082, -0.8357415758833607, -3.073275580698133, 1.045012784247494),
    @(1000, "struggle", 2.214876651763916, 2.617048025131226, -4.977540016174316, -0.3463418929248446, -2.08665958715945, 0.5867890514889552),
    @(1100, "struggle", -4.946963310241699, -1.464599609375, 0.2086009979248047, -0.1247474120891824, -1.187129327229091, 0.3211249946331491),
    @(1200, "struggle", -1.311498880386352, -0.6331030130386353, -0.9298524856567384, -0.08147882624548454, -0.2940848761675368, 0.3946158166746729),
    @(1300, "struggle", -0.7284626960754395, 1.614851474761963, 4.958660125732422, 0.217271684383857, 1.024243365745133, -0.1530809580367363),
    @(1400, "struggle", -11.12845420837402, 3.947137832641602, 8.631446838378906, 0.7111998541014513, 5.599125376769453, -1.684482739440025),
    @(1500, "struggle", 5.936349868774414, -0.3100541830062866, -5.833556175231934, -0.03431442684056363, 5.303968478222288, -0.2843577630665928),
    @(1600, "struggle", -1.618814945220947, -2.745197772979736, -3.147943019866944, 0.09735821826116758, 3.700285123319055, 0.6092508264950334),
    @(1700, "struggle", -3.717551708221436, -0.5804144740104675, 2.004154682159424, 0.4117356596552622, 2.450753956424948, 0.3930325690580906),
    @(1800, "struggle", -2.254782676696777, 2.873416900634766, -8.692718505859375, -0.05785143649091724, 1.907227247953421, 0.5867890514889552),
    @(1900, "struggle", -8.337739944458008, -3.719696521759033, -4.337332725524902, 0.4829513220762701, 0.2850683629512847, -0.8451413950141577),
    @(2000, "struggle", 1.585304737091064, 1.517187714576721, 7.499239444732666, 0.1798686522002135, -4.301874998880876, -1.544167586735314),
    @(2100, "struggle", -2.441369771957397, 1.091274261474609, 0.2500591278076172, -0.1045452207326868, -6.587047031947534, -0.4474151475089096),
    @(2200, "struggle", -2.391006469726562, 4.326517105102539, 1.920782089233398, -0.5947178358934316, -5.224094493048558, 0.5132795298586089),
    @(2300, "struggle", -8.625715255737305, 1.131956696510315, 10.77142906188965, -0.8548342390936281, -2.185787964840323, 0.3064268468594082),
    @(2400, "struggle", -2.745645046234131, -3.043209075927734, -2.123822450637817, 0.3455876613758052, -2.538163804582192, -0.3575649772371564),
    @(2500, "struggle", 0.9841413497924804, 2.591554641723633, -6.520083427429199, 0.1558454343859045, -0.2726017300571888, -0.3659176315580095),
    @(2600, "struggle", -5.384963035583496, 1.437288880348206, -2.568338394165039, 0.8777634601811899, 3.256778053787261, -1.151989024512615),
    @(2700, "struggle", 0.6416101455688477, 0.1077315807342529, -3.572567462921143, 1.138076220239913, 4.945667840996566, -0.5911586260309147),
    @(2800, "struggle", 1.180892944335938, -0.3624088764190674, 1.944910764694214, 0.5868015289306701, 4.023616756711703, 0.7568838426044971),
    @(2900, "struggle", -0.6099348068237305, -0.0995303392410278, 1.559979677200317, -0.3439888250538894, 1.417871174155451, 1.152574896812441)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Host "UsedRange: $($ws.UsedRange.Address())"